# Scheduled runner: refresh market-board derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) on the per-job Profits sheets with the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 138.88235
$ws.Range("I39").Value = 74
$ws.Range("J39").Value = 231.57143
$ws.Range("K39").Value = 222
$ws.Range("L39").Value = 694.71429
$ws.Range("M39").Value = 74
$ws.Range("N39").Value = -1286.71429
$ws.Range("H43").Value = 8078.6875
$ws.Range("J43").Value = 13357
$ws.Range("L43").Value = 13357
$ws.Range("N43").Value = -13495
$ws.Range("H69").Value = 4503
$ws.Range("I69").Value = 6250
$ws.Range("J69").Value = 3338.3333
$ws.Range("K69").Value = 18750
$ws.Range("L69").Value = 10014.9999
$ws.Range("M69").Value = -17876
$ws.Range("N69").Value = -11762.9999
$ws.Range("H72").Value = 4503
$ws.Range("I72").Value = 6250
$ws.Range("J72").Value = 3338.3333
$ws.Range("K72").Value = 56250
$ws.Range("L72").Value = 30044.9997
$ws.Range("M72").Value = -51882
$ws.Range("N72").Value = -38780.9997
$ws.Range("H98").Value = 1550.6154
$ws.Range("I98").Value = 1596.5
$ws.Range("K98").Value = 1596.5
$ws.Range("M98").Value = -98.5
$ws.Range("H107").Value = 4887.8857
$ws.Range("I107").Value = 6391.04
$ws.Range("K107").Value = 6391.04
$ws.Range("M107").Value = -4471.04
$ws.Range("H116").Value = 4600.385
$ws.Range("I116").Value = 5173.1816
$ws.Range("J116").Value = 1450
$ws.Range("K116").Value = 5173.1816
$ws.Range("L116").Value = 1450
$ws.Range("M116").Value = -1731.1816
$ws.Range("N116").Value = -8334
$ws.Range("H122").Value = 1550.6154
$ws.Range("I122").Value = 1596.5
$ws.Range("K122").Value = 4789.5
$ws.Range("M122").Value = -2339.5
$ws.Range("H135").Value = 2740.1052
$ws.Range("I135").Value = 2316.375
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 20847.375
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -18312.375
$ws.Range("N135").Value = -50070
$ws.Range("H138").Value = 4809444
$ws.Range("J138").Value = 8930690
$ws.Range("L138").Value = 26792070
$ws.Range("N138").Value = -26802350
$ws.Range("H139").Value = 37950
$ws.Range("J139").Value = 37950
$ws.Range("L139").Value = 37950
$ws.Range("N139").Value = -48230
$ws.Range("H140").Value = 79060
$ws.Range("J140").Value = 79060
$ws.Range("L140").Value = 79060
$ws.Range("N140").Value = -89420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 962.3333
$ws.Range("I2").Value = 850.8333
$ws.Range("J2").Value = 1185.3334
$ws.Range("K2").Value = 850.8333
$ws.Range("L2").Value = 1185.3334
$ws.Range("M2").Value = -737.8333
$ws.Range("N2").Value = -1411.3334
$ws.Range("H61").Value = 41672084
$ws.Range("I61").Value = 62505628
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 62505628
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -62505416
$ws.Range("N61").Value = -5424
$ws.Range("H116").Value = 962.3333
$ws.Range("I116").Value = 850.8333
$ws.Range("J116").Value = 1185.3334
$ws.Range("K116").Value = 850.8333
$ws.Range("L116").Value = 1185.3334
$ws.Range("M116").Value = 1443.1667
$ws.Range("N116").Value = -5773.3334
$ws.Range("H122").Value = 10759.647
$ws.Range("I122").Value = 11135.714
$ws.Range("K122").Value = 33407.142
$ws.Range("M122").Value = -30957.142
$ws.Range("H132").Value = 16669903
$ws.Range("I132").Value = 27780298
$ws.Range("J132").Value = 4310.3335
$ws.Range("K132").Value = 83340894
$ws.Range("L132").Value = 12931.0005
$ws.Range("M132").Value = -83338364
$ws.Range("N132").Value = -17991.0005
$ws.Range("H136").Value = 41672084
$ws.Range("I136").Value = 62505628
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 187516884
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -187514334
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 962.3333
$ws.Range("I3").Value = 850.8333
$ws.Range("J3").Value = 1185.3334
$ws.Range("K3").Value = 850.8333
$ws.Range("L3").Value = 1185.3334
$ws.Range("M3").Value = -736.8333
$ws.Range("N3").Value = -1413.3334
$ws.Range("H105").Value = 4319.946
$ws.Range("I105").Value = 3253.9
$ws.Range("J105").Value = 4714.778
$ws.Range("K105").Value = 3253.9
$ws.Range("L105").Value = 4714.778
$ws.Range("M105").Value = -1506.9
$ws.Range("N105").Value = -8208.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23820572
$ws.Range("I31").Value = 12876.637
$ws.Range("K31").Value = 12876.637
$ws.Range("M31").Value = -12581.637
$ws.Range("H34").Value = 23820572
$ws.Range("I34").Value = 12876.637
$ws.Range("K34").Value = 12876.637
$ws.Range("M34").Value = -12674.637
$ws.Range("H105").Value = 2194
$ws.Range("I105").Value = 1060
$ws.Range("J105").Value = 2950
$ws.Range("K105").Value = 1060
$ws.Range("L105").Value = 2950
$ws.Range("M105").Value = 687
$ws.Range("N105").Value = -6444
$ws.Range("H134").Value = 2081.6316
$ws.Range("I134").Value = 1815.6875
$ws.Range("K134").Value = 5447.0625
$ws.Range("M134").Value = -2912.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 729.05884
$ws.Range("I107").Value = 199
$ws.Range("J107").Value = 762.1875
$ws.Range("K107").Value = 597
$ws.Range("L107").Value = 2286.5625
$ws.Range("M107").Value = 1323
$ws.Range("N107").Value = -6126.5625
$ws.Range("H131").Value = 920
$ws.Range("I131").Value = 732
$ws.Range("K131").Value = 2196
$ws.Range("M131").Value = 2844
$ws.Range("H132").Value = 739.8570999999999
$ws.Range("I132").Value = 435.8
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 3922.2
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -1392.2
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9003630
$ws.Range("I80").Value = 17546288
$ws.Range("K80").Value = 17546288
$ws.Range("M80").Value = -17545290
$ws.Range("H83").Value = 9003630
$ws.Range("I83").Value = 17546288
$ws.Range("K83").Value = 87731440
$ws.Range("M83").Value = -87726448
$ws.Range("H97").Value = 884.4828
$ws.Range("I97").Value = 883.2857
$ws.Range("K97").Value = 883.2857
$ws.Range("M97").Value = -387.2857
$ws.Range("H102").Value = 4100
$ws.Range("I102").Value = 5206
$ws.Range("J102").Value = 1335
$ws.Range("K102").Value = 5206
$ws.Range("L102").Value = 1335
$ws.Range("M102").Value = -3584
$ws.Range("N102").Value = -4579
$ws.Range("H122").Value = 3176608.2
$ws.Range("I122").Value = 5129988
$ws.Range("J122").Value = 2366.125
$ws.Range("K122").Value = 15389964
$ws.Range("L122").Value = 7098.375
$ws.Range("M122").Value = -15387514
$ws.Range("N122").Value = -11998.375
$ws.Range("H132").Value = 5509
$ws.Range("I132").Value = 4469.8
$ws.Range("J132").Value = 6251.2856
$ws.Range("K132").Value = 13409.4
$ws.Range("L132").Value = 18753.8568
$ws.Range("M132").Value = -10879.4
$ws.Range("N132").Value = -23813.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6984.0713
$ws.Range("J40").Value = 9666.666999999999
$ws.Range("L40").Value = 9666.666999999999
$ws.Range("N40").Value = -9938.666999999999
$ws.Range("H82").Value = 2176.647
$ws.Range("I82").Value = 1974.909
$ws.Range("K82").Value = 1974.909
$ws.Range("M82").Value = -1613.909
$ws.Range("H85").Value = 2176.647
$ws.Range("I85").Value = 1974.909
$ws.Range("K85").Value = 1974.909
$ws.Range("M85").Value = -726.9090000000001
$ws.Range("H139").Value = 39409.4
$ws.Range("J139").Value = 39271.555
$ws.Range("L139").Value = 39271.555
$ws.Range("N139").Value = -49551.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3383.7693
$ws.Range("I122").Value = 3320.4
$ws.Range("K122").Value = 9961.200000000001
$ws.Range("M122").Value = -7511.200000000001
